$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their text representation
# (values like "1.002" or "29.373.73" must not be reinterpreted as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.373.73'
$ws.Range('D3').Value = '1.843.37'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '240.52'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = '0.6266'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.07473'
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').Value = '0.2888'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '0.07725'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '1.843.51'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').Value = '4.986'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = '0.6780'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').Value = '0.00001025'
$ws.Range('E15').Value = '  -5.46%  '
$ws.Range('D16').Value = '82.16'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').Value = '2.108.80'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '6.082'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '29.401.01'
$ws.Range('D20').Value = '228.66'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '12.26'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '7.363'
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '158.92'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '0.1380'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '8.378'
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').Value = '17.54'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '1.395'
$ws.Range('E29').Value = '  +3.05%  '
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').Value = '0.05680'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').Value = '4.093'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').Value = '4.041'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').Value = '1.818'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('D36').Value = '0.6925'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('D37').Value = '2.589'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '2.853'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('D39').Value = '1.251.98'
$ws.Range('E39').Value = '  +2.13%  '
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('D41').Value = '6.509'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').Value = '0.9030'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '2.009.54'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = '101.16'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '65.59'
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').Value = '7.068'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').Value = '0.1158'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000115'
$ws.Range('E49').Value = '  -3.69%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.932'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = '0.3928'
$ws.Range('E51').Value = '  -2.35%  '
